# Update column G ("K") values on Sheet1 for rows 2-23.
# New values regenerated after switching from "Strike#" to "K" calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 12
    3  = 3
    4  = 6
    5  = 8
    6  = 10
    7  = 3
    8  = 8
    9  = 3
    10 = 7
    11 = 5
    12 = 7
    13 = 3
    14 = 4
    15 = 4
    16 = 7
    17 = 12
    18 = 7
    19 = 4
    20 = 6
    21 = 5
    22 = 6
    23 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
